# Auto-applies the cryptos.xlsx price/volume/ranking update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($addr, $val) {
    # Writes $val into $addr as literal text, even when it parses as a number
    # (e.g. "114.64"), then restores the cell to its original (default) style so
    # no stray number-format / style index is left behind on the cell.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

function Set-Text($addr, $val) {
    # Plain strings (urls, names, percent strings) that Excel will not
    # mis-parse as numbers/dates can be written directly.
    $ws.Range($addr).Value = $val
}

Set-Text 'D2' '43.637.21'
Set-Text 'E2' '  -0.07%  '
Set-Text 'D3' '2.293.29'
Set-Text 'E3' '  +0.36%  '
Set-Text 'E4' '  +0.11%  '
Set-PlainText 'D5' '114.64'
Set-Text 'E5' '  +19.20%  '
Set-PlainText 'D6' '268.91'
Set-Text 'E6' '  +0.81%  '
Set-PlainText 'D7' '0.616'
Set-Text 'E8' '  +0.18%  '
Set-PlainText 'D10' '48.56'
Set-Text 'E10' '  +6.49%  '
Set-PlainText 'D11' '0.0939'
Set-Text 'E11' '  +0.50%  '
Set-PlainText 'D12' '8.65'
Set-Text 'E12' '  +10.88%  '
Set-Text 'E13' '  +1.08%  '
Set-PlainText 'D14' '15.58'
Set-Text 'E14' '  +3.15%  '
Set-Text 'D15' '2.637.74'
Set-Text 'E15' '  +0.40%  '
Set-PlainText 'D16' '0.847'
Set-Text 'E16' '  +0.63%  '
Set-Text 'D17' '2.292.85'
Set-Text 'E17' '  +0.11%  '
Set-Text 'D18' '43.645.38'
Set-Text 'E18' '  +0.11%  '
Set-Text 'E19' '  +2.30%  '
Set-PlainText 'D20' '6.55'
Set-Text 'E20' '  +5.56%  '
Set-PlainText 'D21' '72.55'
Set-Text 'E21' '  +0.81%  '
Set-PlainText 'D22' '2.50'
Set-Text 'E22' '  +2.29%  '
Set-PlainText 'D23' '233.24'
Set-Text 'E23' '  +0.19%  '
Set-PlainText 'D24' '9.61'
Set-Text 'E24' '  +5.06%  '
Set-PlainText 'D25' '2.82'
Set-Text 'E25' '  +13.21%  '
Set-PlainText 'D26' '0.999'
Set-Text 'E26' '  +0.01%  '
Set-PlainText 'D27' '11.56'
Set-Text 'E27' '  +3.97%  '
Set-PlainText 'D28' '42.34'
Set-Text 'E28' '  +5.47%  '
Set-PlainText 'D29' '3.40'
Set-PlainText 'D30' '2.28'
Set-Text 'E30' '  +0.51%  '
Set-PlainText 'D31' '176.28'
Set-Text 'E31' '  -0.03%  '
Set-PlainText 'D32' '0.0936'
Set-Text 'E32' '  +4.84%  '
Set-PlainText 'D33' '21.60'
Set-Text 'E33' '  -1.03%  '
Set-PlainText 'D34' '5.53'
Set-Text 'E34' '  +3.36%  '
Set-Text 'E35' '  +0.74%  '
Set-Text 'E36' '  +9.81%  '
Set-Text 'E37' '  +0.89%  '
Set-PlainText 'D38' '0.0357'
Set-Text 'E38' '  +0.89%  '
Set-PlainText 'D39' '3.85'
Set-Text 'E39' '  +13.34%  '
Set-Text 'B40' 'LidoDAOToken'
Set-Text 'C40' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-PlainText 'D40' '2.42'
Set-Text 'E40' '  +5.05%  '
Set-Text 'B41' 'Algorand'
Set-Text 'C41' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-PlainText 'D41' '0.243'
Set-Text 'E41' '  +3.59%  '
Set-PlainText 'D42' '13.78'
Set-Text 'E42' '  +12.52%  '
Set-Text 'B43' 'MultiversX'
Set-Text 'C43' 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-PlainText 'D43' '73.70'
Set-Text 'E43' '  +14.79%  '
Set-PlainText 'D44' '1.45'
Set-Text 'E44' '  +8.05%  '
Set-PlainText 'D45' '6.00'
Set-Text 'E45' '  +15.42%  '
Set-Text 'E46' '  +0.10%  '
Set-PlainText 'D47' '8.76'
Set-Text 'E47' '  -0.04%  '
Set-Text 'B48' 'Cronos'
Set-Text 'C48' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-PlainText 'D48' '0.100'
Set-Text 'E48' '  -1.70%  '
Set-Text 'B49' 'Aave'
Set-Text 'C49' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-PlainText 'D49' '102.11'
Set-Text 'E49' '  +5.47%  '
Set-Text 'E50' '  +3.77%  '
Set-PlainText 'D51' '0.451'
Set-Text 'E51' '  +5.08%  '
